$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing data occupies rows 2-101 (header in row 1).
# We append 45 new rows (102-146) following the same repeating pattern
# seen in the source data: column A cycles through 9 id values,
# column B increments by 1 each row, and C-G repeat constant values.

$idCycle = @(10002, 10003, 10004, 10005, 10006, 10007, 10008, 10009, 10010)

$startRow = 102
$endRow = 146
$startDeviceId = 3000121

for ($r = $startRow; $r -le $endRow; $r++) {
    $cycleIndex = ($r - $startRow) % 9
    $regCenterId = $idCycle[$cycleIndex]
    $deviceId = $startDeviceId + ($r - $startRow)

    $ws.Cells.Item($r, 1).Value = $regCenterId
    $ws.Cells.Item($r, 2).Value = $deviceId
    $ws.Cells.Item($r, 3).Value = "eng"
    $ws.Cells.Item($r, 4).Value = $true
    $ws.Cells.Item($r, 5).Value = "superadmin"
    $ws.Cells.Item($r, 6).Value = "now()"
    $ws.Cells.Item($r, 7).Value = "now()"
}

# Reflect the final selection/scroll position left behind in the sheet view.
$ws.Range("A102:B146").Select()
$excel.ActiveWindow.ScrollRow = 128
